# Corrected excel sheets for application fix issues
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet: move active selection from C4 to D3
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("D3").Select()

# ---------------------------------------------------------------------------
# Repayment schedule sheet: add the missing "O" column figures and move the
# active selection from E6 to the whole of row 15
# ---------------------------------------------------------------------------
$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Row 2 and 3 only need a (blank) cell added in column P, matching the
# existing blank formatting used across the row (copy format from column N,
# which already carries the correct style for these two header/blank rows).
foreach ($r in 2, 3) {
    $src = $wsRepay.Cells.Item($r, 14)
    $dst = $wsRepay.Cells.Item($r, 16)
    $src.Copy($dst)
}

# Rows 4 through 14 are missing column O (the sheet jumps straight from N to
# P) - fill it in with the same "0" figure used throughout the schedule,
# copying the formatting from the neighbouring N column so the new cells
# match the existing style.
for ($r = 4; $r -le 14; $r++) {
    $src = $wsRepay.Cells.Item($r, 14)
    $dst = $wsRepay.Cells.Item($r, 15)
    $src.Copy($dst)
    $dst.Value = 0
}

$wsRepay.Rows.Item(15).Select()

# ---------------------------------------------------------------------------
# Transactions sheet: renumber the transaction ids and move the active
# selection from D4 to D3
# ---------------------------------------------------------------------------
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("A2").Value = 7078
$wsTransactions.Range("A3").Value = 7077
$wsTransactions.Range("A4").Value = 7076
$wsTransactions.Range("D3").Select()

# Restore Transactions as the active/visible sheet (it was the tab selected
# before these edits were made).
$wsTransactions.Activate()
